# Auto-generated edit script applying scheduled runner data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1950.1029
$ws.Range("I15").Value = 1950.1029
$ws.Range("K15").Value = 5850.3087
$ws.Range("M15").Value = -5681.3087

$ws.Range("H17").Value = 1473532.1
$ws.Range("J17").Value = 1541530.1
$ws.Range("L17").Value = 4624590.300000001
$ws.Range("N17").Value = -4624926.300000001

$ws.Range("H33").Value = 77.125
$ws.Range("I33").Value = 77.125
$ws.Range("K33").Value = 77.125
$ws.Range("M33").Value = 151.875

$ws.Range("H40").Value = 1673.8462
$ws.Range("I40").Value = 1320
$ws.Range("J40").Value = 2240
$ws.Range("K40").Value = 1320
$ws.Range("L40").Value = 2240
$ws.Range("M40").Value = -1145
$ws.Range("N40").Value = -2590

$ws.Range("H62").Value = 2865.9375
$ws.Range("I62").Value = 2200.5557
$ws.Range("J62").Value = 3721.4285
$ws.Range("K62").Value = 2200.5557
$ws.Range("L62").Value = 3721.4285
$ws.Range("M62").Value = -1576.5557
$ws.Range("N62").Value = -4969.4285

$ws.Range("H65").Value = 2865.9375
$ws.Range("I65").Value = 2200.5557
$ws.Range("J65").Value = 3721.4285
$ws.Range("K65").Value = 11002.7785
$ws.Range("L65").Value = 18607.1425
$ws.Range("M65").Value = -7882.7785
$ws.Range("N65").Value = -24847.1425

$ws.Range("H70").Value = 1750
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1750
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5250
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -5790

$ws.Range("H73").Value = 1750
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1750
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5250
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -7122

$ws.Range("H88").Value = 1429.8235
$ws.Range("I88").Value = 1164.5714
$ws.Range("J88").Value = 1615.5
$ws.Range("K88").Value = 1164.5714
$ws.Range("L88").Value = 1615.5
$ws.Range("M88").Value = -758.5714
$ws.Range("N88").Value = -2427.5

$ws.Range("H91").Value = 1429.8235
$ws.Range("I91").Value = 1164.5714
$ws.Range("J91").Value = 1615.5
$ws.Range("K91").Value = 1164.5714
$ws.Range("L91").Value = 1615.5
$ws.Range("M91").Value = 239.4286
$ws.Range("N91").Value = -4423.5

$ws.Range("H96").Value = 35715224
$ws.Range("I96").Value = 41667504
$ws.Range("K96").Value = 125002512
$ws.Range("M96").Value = -125001139

$ws.Range("H97").Value = 2481.375
$ws.Range("J97").Value = 2481.375
$ws.Range("L97").Value = 7444.125
$ws.Range("N97").Value = -8436.125

$ws.Range("H137").Value = 48034.773
$ws.Range("I137").Value = 2854.5715
$ws.Range("K137").Value = 8563.7145
$ws.Range("M137").Value = -6013.7145

$ws.Range("H138").Value = 1396.6202
$ws.Range("I138").Value = 515.4792
$ws.Range("K138").Value = 1546.4376
$ws.Range("M138").Value = 3593.5624

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1305.6364
$ws.Range("I2").Value = 1120.96
$ws.Range("K2").Value = 1120.96
$ws.Range("M2").Value = -1007.96

$ws.Range("H5").Value = 98.333336
$ws.Range("I5").Value = 98.333336
$ws.Range("K5").Value = 98.333336
$ws.Range("M5").Value = 13.666664

$ws.Range("H32").Value = 18929.934
$ws.Range("I32").Value = 19954.219
$ws.Range("J32").Value = 7662.8
$ws.Range("K32").Value = 19954.219
$ws.Range("L32").Value = 7662.8
$ws.Range("M32").Value = -19667.219
$ws.Range("N32").Value = -8236.799999999999

$ws.Range("H45").Value = 2490.75
$ws.Range("I45").Value = 2736.4546
$ws.Range("K45").Value = 2736.4546
$ws.Range("M45").Value = -2359.4546

$ws.Range("H116").Value = 1305.6364
$ws.Range("I116").Value = 1120.96
$ws.Range("K116").Value = 1120.96
$ws.Range("M116").Value = 1173.04

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1305.6364
$ws.Range("I3").Value = 1120.96
$ws.Range("K3").Value = 1120.96
$ws.Range("M3").Value = -1006.96

$ws.Range("H4").Value = 98.333336
$ws.Range("I4").Value = 98.333336
$ws.Range("K4").Value = 98.333336
$ws.Range("M4").Value = 16.666664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9371.413
$ws.Range("I31").Value = 12997.186
$ws.Range("K31").Value = 12997.186
$ws.Range("M31").Value = -12702.186

$ws.Range("H34").Value = 9371.413
$ws.Range("I34").Value = 12997.186
$ws.Range("K34").Value = 12997.186
$ws.Range("M34").Value = -12795.186

$ws.Range("H58").Value = 11542.426
$ws.Range("I58").Value = 935.67645
$ws.Range("J58").Value = 39283.152
$ws.Range("K58").Value = 935.67645
$ws.Range("L58").Value = 39283.152
$ws.Range("M58").Value = -732.67645
$ws.Range("N58").Value = -39689.152

$ws.Range("H86").Value = 6180540.5
$ws.Range("I86").Value = 2690.4
$ws.Range("J86").Value = 13902853
$ws.Range("K86").Value = 2690.4
$ws.Range("L86").Value = 13902853
$ws.Range("M86").Value = -1567.4
$ws.Range("N86").Value = -13905099

$ws.Range("H89").Value = 6180540.5
$ws.Range("I89").Value = 2690.4
$ws.Range("J89").Value = 13902853
$ws.Range("K89").Value = 13452
$ws.Range("L89").Value = 69514265
$ws.Range("M89").Value = -7836
$ws.Range("N89").Value = -69525497

$ws.Range("H134").Value = 1297.4878
$ws.Range("I134").Value = 974.15
$ws.Range("K134").Value = 2922.45
$ws.Range("M134").Value = -387.4499999999998

$ws.Range("H136").Value = 11542.426
$ws.Range("I136").Value = 935.67645
$ws.Range("J136").Value = 39283.152
$ws.Range("K136").Value = 2807.02935
$ws.Range("L136").Value = 117849.456
$ws.Range("M136").Value = -257.0293500000002
$ws.Range("N136").Value = -122949.456

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 92106.2
$ws.Range("J36").Value = 151901.33
$ws.Range("L36").Value = 455703.99
$ws.Range("N36").Value = -456041.99

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H98").Value = 448.33334
$ws.Range("I98").Value = 256.75
$ws.Range("J98").Value = 601.6
$ws.Range("K98").Value = 770.25
$ws.Range("L98").Value = 1804.8
$ws.Range("M98").Value = 727.75
$ws.Range("N98").Value = -4800.8

$ws.Range("H107").Value = 20250
$ws.Range("J107").Value = 250
$ws.Range("L107").Value = 750
$ws.Range("N107").Value = -4590

$ws.Range("H131").Value = 747.9596
$ws.Range("I131").Value = 574.75
$ws.Range("J131").Value = 755.2526
$ws.Range("K131").Value = 1724.25
$ws.Range("L131").Value = 2265.7578
$ws.Range("M131").Value = 3315.75
$ws.Range("N131").Value = -12345.7578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 29414294
$ws.Range("I102").Value = 33336094
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 33336094
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = -33334472
$ws.Range("N102").Value = -4044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 29250
$ws.Range("J59").Value = 29250
$ws.Range("L59").Value = 29250
$ws.Range("N59").Value = -30558

$ws.Range("H61").Value = 5591.9165
$ws.Range("I61").Value = 3585.7144
$ws.Range("K61").Value = 3585.7144
$ws.Range("M61").Value = -3383.7144

$ws.Range("H113").Value = 5591.9165
$ws.Range("I113").Value = 3585.7144
$ws.Range("K113").Value = 3585.7144
$ws.Range("M113").Value = -1415.7144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 33084
$ws.Range("J27").Value = 33084
$ws.Range("L27").Value = 33084
$ws.Range("N27").Value = -33222

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
